$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data (row 210) spanning columns A (1) through DX (128),
# matching the column order of the existing header row.
$rowValues = @(44104, 829679, 2744, 114694, 67554, 269142, 28879, 7137, 5947, 8488, 9490, 20013, 3997, 23953, 33306, 8136, 11492, 15339, 15258, 18420, 15796, 3790, 3708, 11002, 31229, 14081, 12384, 62704, 2439, 1478, 769, 477, 874, 529, 774, 2071, 5921, 38163, 10181, 2574, 48345, 1164, 23018, 1535, 10452, 1683, 1615, 8613, 2051, 965, 2503, 2695, 66318, 14214, 6849, 9984, 7403, 257, 1465, 2737, 745, 2169, 10001, 9601, 10761, 14381, 1972, 908, 14241, 11666, 13633, 3274, 2305, 6198, 5025, 2632, 6230, 3882, 2285, 1270, 3080, 2252, 2109, 1963, 6615, 2239, 1519, 1850, 2151, 2308, 2682, 1838, 1231, 1231, 1108, 3479, 1513, 980, 1156, 1761, 1637, 822, 918, 1385, 1819, 1694, 1671, 1275, 334, 373, 840, 786, 502, 544, 389, 686, 761, 530, 507, 374, 527, 139921, 350353, 20930, 152718, 94611, 47074, 13152)

$row = 210
for ($col = 1; $col -le $rowValues.Length; $col++) {
    $ws.Cells.Item($row, $col).Value = $rowValues[$col - 1]
}

# Move the active selection to the new last row (A210), matching Excel's
# behaviour of landing on the newly entered row after data entry.
$ws.Range("A210").Select()
